$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 8 ("Brasil"), shifting it down to row 9
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the Sergipe data
$ws.Range("A8").Value = "Sergipe"
$ws.Range("B8").Value = "Diferença 2021-2012"
$ws.Range("C8").Value = 11.71108875883282
$ws.Range("D8").Value = "12º"
